$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The NATMI TPM data was regenerated. The table now includes "Resolving-Mac"
# as a possible target cluster (previously it was only a sending cluster),
# producing a full 4x4 sending/target cluster cross product (16 data rows)
# instead of the previous 4x3 cross product (12 data rows). Values in most
# numeric columns changed accordingly.

$numRows = 16
$numCols = 20
$arr = New-Object 'object[,]' $numRows,$numCols

$arr[0,0] = "ECs"
$arr[0,1] = "Ccl2"
$arr[0,2] = "Ackr4"
$arr[0,3] = "ECs"
$arr[0,4] = 3
$arr[0,5] = 1
$arr[0,6] = 3.719288
$arr[0,7] = 11.157864
$arr[0,8] = 0.04235839908674209
$arr[0,9] = 0.04235839908674209
$arr[0,10] = 1
$arr[0,11] = 0.3333333333333333
$arr[0,12] = 0.1143813333333333
$arr[0,13] = 0.343144
$arr[0,14] = 0.128300337591142
$arr[0,15] = 0.1283003375911419
$arr[0,16] = 0.4254171204906667
$arr[0,17] = 3.828754084416
$arr[0,18] = 0.005434596902649329
$arr[0,19] = 0.005434596902649328
$arr[1,0] = "ECs"
$arr[1,1] = "Ccl2"
$arr[1,2] = "Ackr4"
$arr[1,3] = "FAPs"
$arr[1,4] = 3
$arr[1,5] = 1
$arr[1,6] = 3.719288
$arr[1,7] = 11.157864
$arr[1,8] = 0.04235839908674209
$arr[1,9] = 0.04235839908674209
$arr[1,10] = 3
$arr[1,11] = 1
$arr[1,12] = 0.7200953333333334
$arr[1,13] = 2.160286
$arr[1,14] = 0.8077233554817153
$arr[1,15] = 0.8077233554817151
$arr[1,16] = 2.678241932122667
$arr[1,17] = 24.104177389104
$arr[1,18] = 0.03421386824317695
$arr[1,19] = 0.03421386824317694
$arr[2,0] = "ECs"
$arr[2,1] = "Ccl2"
$arr[2,2] = "Ackr4"
$arr[2,3] = "MuSCs"
$arr[2,4] = 3
$arr[2,5] = 1
$arr[2,6] = 3.719288
$arr[2,7] = 11.157864
$arr[2,8] = 0.04235839908674209
$arr[2,9] = 0.04235839908674209
$arr[2,10] = 1
$arr[2,11] = 0.3333333333333333
$arr[2,12] = 0.03357866666666667
$arr[2,13] = 0.100736
$arr[2,14] = 0.03766483694187069
$arr[2,15] = 0.03766483694187069
$arr[2,16] = 0.1248887319893333
$arr[2,17] = 1.123998587904
$arr[2,18] = 0.001595422194720825
$arr[2,19] = 0.001595422194720825
$arr[3,0] = "ECs"
$arr[3,1] = "Ccl2"
$arr[3,2] = "Ackr4"
$arr[3,3] = "Resolving-Mac"
$arr[3,4] = 3
$arr[3,5] = 1
$arr[3,6] = 3.719288
$arr[3,7] = 11.157864
$arr[3,8] = 0.04235839908674209
$arr[3,9] = 0.04235839908674209
$arr[3,10] = 1
$arr[3,11] = 0.3333333333333333
$arr[3,12] = 0.023457
$arr[3,13] = 0.070371
$arr[3,14] = 0.02631146998527222
$arr[3,15] = 0.02631146998527222
$arr[3,16] = 0.087243338616
$arr[3,17] = 0.785190047544
$arr[3,18] = 0.001114511746194997
$arr[3,19] = 0.001114511746194997
$arr[4,0] = "FAPs"
$arr[4,1] = "Ccl2"
$arr[4,2] = "Ackr4"
$arr[4,3] = "ECs"
$arr[4,4] = 3
$arr[4,5] = 1
$arr[4,6] = 26.15942766666667
$arr[4,7] = 78.478283
$arr[4,8] = 0.2979256989470644
$arr[4,9] = 0.2979256989470644
$arr[4,10] = 1
$arr[4,11] = 0.3333333333333333
$arr[4,12] = 0.1143813333333333
$arr[4,13] = 0.343144
$arr[4,14] = 0.128300337591142
$arr[4,15] = 0.1283003375911419
$arr[4,16] = 2.992150215750223
$arr[4,17] = 26.929351941752
$arr[4,18] = 0.03822396775198528
$arr[4,19] = 0.03822396775198528
$arr[5,0] = "FAPs"
$arr[5,1] = "Ccl2"
$arr[5,2] = "Ackr4"
$arr[5,3] = "FAPs"
$arr[5,4] = 3
$arr[5,5] = 1
$arr[5,6] = 26.15942766666667
$arr[5,7] = 78.478283
$arr[5,8] = 0.2979256989470644
$arr[5,9] = 0.2979256989470644
$arr[5,10] = 3
$arr[5,11] = 1
$arr[5,12] = 0.7200953333333334
$arr[5,13] = 2.160286
$arr[5,14] = 0.8077233554817153
$arr[5,15] = 0.8077233554817151
$arr[5,16] = 18.83728178543756
$arr[5,17] = 169.535536068938
$arr[5,18] = 0.2406415452377582
$arr[5,19] = 0.2406415452377581
$arr[6,0] = "FAPs"
$arr[6,1] = "Ccl2"
$arr[6,2] = "Ackr4"
$arr[6,3] = "MuSCs"
$arr[6,4] = 3
$arr[6,5] = 1
$arr[6,6] = 26.15942766666667
$arr[6,7] = 78.478283
$arr[6,8] = 0.2979256989470644
$arr[6,9] = 0.2979256989470644
$arr[6,10] = 1
$arr[6,11] = 0.3333333333333333
$arr[6,12] = 0.03357866666666667
$arr[6,13] = 0.100736
$arr[6,14] = 0.03766483694187069
$arr[6,15] = 0.03766483694187069
$arr[6,16] = 0.8783987018097779
$arr[6,17] = 7.905588316288001
$arr[6,18] = 0.01122132287163403
$arr[6,19] = 0.01122132287163403
$arr[7,0] = "FAPs"
$arr[7,1] = "Ccl2"
$arr[7,2] = "Ackr4"
$arr[7,3] = "Resolving-Mac"
$arr[7,4] = 3
$arr[7,5] = 1
$arr[7,6] = 26.15942766666667
$arr[7,7] = 78.478283
$arr[7,8] = 0.2979256989470644
$arr[7,9] = 0.2979256989470644
$arr[7,10] = 1
$arr[7,11] = 0.3333333333333333
$arr[7,12] = 0.023457
$arr[7,13] = 0.070371
$arr[7,14] = 0.02631146998527222
$arr[7,15] = 0.02631146998527222
$arr[7,16] = 0.6136216947770001
$arr[7,17] = 5.522595252993001
$arr[7,18] = 0.007838863085686931
$arr[7,19] = 0.007838863085686931
$arr[8,0] = "MuSCs"
$arr[8,1] = "Ccl2"
$arr[8,2] = "Ackr4"
$arr[8,3] = "ECs"
$arr[8,4] = 3
$arr[8,5] = 1
$arr[8,6] = 6.299630666666666
$arr[8,7] = 18.898892
$arr[8,8] = 0.07174552491706633
$arr[8,9] = 0.07174552491706633
$arr[8,10] = 1
$arr[8,11] = 0.3333333333333333
$arr[8,12] = 0.1143813333333333
$arr[8,13] = 0.343144
$arr[8,14] = 0.128300337591142
$arr[8,15] = 0.1283003375911419
$arr[8,16] = 0.7205601551608889
$arr[8,17] = 6.485041396448
$arr[8,18] = 0.009204975067513298
$arr[8,19] = 0.009204975067513296
$arr[9,0] = "MuSCs"
$arr[9,1] = "Ccl2"
$arr[9,2] = "Ackr4"
$arr[9,3] = "FAPs"
$arr[9,4] = 3
$arr[9,5] = 1
$arr[9,6] = 6.299630666666666
$arr[9,7] = 18.898892
$arr[9,8] = 0.07174552491706633
$arr[9,9] = 0.07174552491706633
$arr[9,10] = 3
$arr[9,11] = 1
$arr[9,12] = 0.7200953333333334
$arr[9,13] = 2.160286
$arr[9,14] = 0.8077233554817153
$arr[9,15] = 0.8077233554817151
$arr[9,16] = 4.536334644790222
$arr[9,17] = 40.827011803112
$arr[9,18] = 0.05795053612680983
$arr[9,19] = 0.05795053612680982
$arr[10,0] = "MuSCs"
$arr[10,1] = "Ccl2"
$arr[10,2] = "Ackr4"
$arr[10,3] = "MuSCs"
$arr[10,4] = 3
$arr[10,5] = 1
$arr[10,6] = 6.299630666666666
$arr[10,7] = 18.898892
$arr[10,8] = 0.07174552491706633
$arr[10,9] = 0.07174552491706633
$arr[10,10] = 1
$arr[10,11] = 0.3333333333333333
$arr[10,12] = 0.03357866666666667
$arr[10,13] = 0.100736
$arr[10,14] = 0.03766483694187069
$arr[10,15] = 0.03766483694187069
$arr[10,16] = 0.2115331982791111
$arr[10,17] = 1.903798784512
$arr[10,18] = 0.002702283497310224
$arr[10,19] = 0.002702283497310224
$arr[11,0] = "MuSCs"
$arr[11,1] = "Ccl2"
$arr[11,2] = "Ackr4"
$arr[11,3] = "Resolving-Mac"
$arr[11,4] = 3
$arr[11,5] = 1
$arr[11,6] = 6.299630666666666
$arr[11,7] = 18.898892
$arr[11,8] = 0.07174552491706633
$arr[11,9] = 0.07174552491706633
$arr[11,10] = 1
$arr[11,11] = 0.3333333333333333
$arr[11,12] = 0.023457
$arr[11,13] = 0.070371
$arr[11,14] = 0.02631146998527222
$arr[11,15] = 0.02631146998527222
$arr[11,16] = 0.147770436548
$arr[11,17] = 1.329933928932
$arr[11,18] = 0.001887730225432991
$arr[11,19] = 0.001887730225432991
$arr[12,0] = "Resolving-Mac"
$arr[12,1] = "Ccl2"
$arr[12,2] = "Ackr4"
$arr[12,3] = "ECs"
$arr[12,4] = 3
$arr[12,5] = 1
$arr[12,6] = 51.62686066666667
$arr[12,7] = 154.880582
$arr[12,8] = 0.5879703770491272
$arr[12,9] = 0.5879703770491272
$arr[12,10] = 1
$arr[12,11] = 0.3333333333333333
$arr[12,12] = 0.1143813333333333
$arr[12,13] = 0.343144
$arr[12,14] = 0.128300337591142
$arr[12,15] = 0.1283003375911419
$arr[12,16] = 5.905149158867555
$arr[12,17] = 53.146342429808
$arr[12,18] = 0.07543679786899404
$arr[12,19] = 0.07543679786899403
$arr[13,0] = "Resolving-Mac"
$arr[13,1] = "Ccl2"
$arr[13,2] = "Ackr4"
$arr[13,3] = "FAPs"
$arr[13,4] = 3
$arr[13,5] = 1
$arr[13,6] = 51.62686066666667
$arr[13,7] = 154.880582
$arr[13,8] = 0.5879703770491272
$arr[13,9] = 0.5879703770491272
$arr[13,10] = 3
$arr[13,11] = 1
$arr[13,12] = 0.7200953333333334
$arr[13,13] = 2.160286
$arr[13,14] = 0.8077233554817153
$arr[13,15] = 0.8077233554817151
$arr[13,16] = 37.17626144071689
$arr[13,17] = 334.586352966452
$arr[13,18] = 0.4749174058739704
$arr[13,19] = 0.4749174058739703
$arr[14,0] = "Resolving-Mac"
$arr[14,1] = "Ccl2"
$arr[14,2] = "Ackr4"
$arr[14,3] = "MuSCs"
$arr[14,4] = 3
$arr[14,5] = 1
$arr[14,6] = 51.62686066666667
$arr[14,7] = 154.880582
$arr[14,8] = 0.5879703770491272
$arr[14,9] = 0.5879703770491272
$arr[14,10] = 1
$arr[14,11] = 0.3333333333333333
$arr[14,12] = 0.03357866666666667
$arr[14,13] = 0.100736
$arr[14,14] = 0.03766483694187069
$arr[14,15] = 0.03766483694187069
$arr[14,16] = 1.733561145372444
$arr[14,17] = 15.602050308352
$arr[14,18] = 0.02214580837820561
$arr[14,19] = 0.02214580837820561
$arr[15,0] = "Resolving-Mac"
$arr[15,1] = "Ccl2"
$arr[15,2] = "Ackr4"
$arr[15,3] = "Resolving-Mac"
$arr[15,4] = 3
$arr[15,5] = 1
$arr[15,6] = 51.62686066666667
$arr[15,7] = 154.880582
$arr[15,8] = 0.5879703770491272
$arr[15,9] = 0.5879703770491272
$arr[15,10] = 1
$arr[15,11] = 0.3333333333333333
$arr[15,12] = 0.023457
$arr[15,13] = 0.070371
$arr[15,14] = 0.02631146998527222
$arr[15,15] = 0.02631146998527222
$arr[15,16] = 1.211011270658
$arr[15,17] = 10.899101435922
$arr[15,18] = 0.0154703649279573
$arr[15,19] = 0.0154703649279573

$startCell = $ws.Cells.Item(2, 1)
$endCell = $ws.Cells.Item(1 + $numRows, $numCols)
$rng = $ws.Range($startCell, $endCell)
$rng.Value2 = $arr
